# Regenerate the "K" column (column G) values in the save_data sheet.
# This mirrors an upstream data regen that swapped "Strike#" for "K" and
# recalculated the statistic written into each row of column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G (header "K"), rows 2-39, in order.
$newK = @(5,9,3,9,7,7,11,4,6,3,8,5,10,11,5,5,7,8,8,6,5,8,2,6,4,5,8,8,4,13,3,8,11,3,7,6,7,4)

$row = 2
foreach ($val in $newK) {
    $ws.Cells.Item($row, 7).Value = $val
    $row++
}
